# The source diff for this change only reorders XML namespace declarations
# and element attributes (e.g. <w:color w:val=".." w:themeColor=".." .../>
# becoming <w:color w:themeColor=".." w:themeShade=".." w:val=".."/>,
# <w:headerReference w:type=".." r:id=".."/> becoming
# <w:headerReference r:id=".." w:type=".."/>, etc.) across document.xml,
# footnotes.xml, header1.xml and styles.xml. The commit message confirms
# this: "Fixed POI packaging and upgraded to POI 3.15" -- i.e. the test
# fixture was re-saved by a newer Apache POI version that serializes
# namespaces/attributes in sorted order, with no change whatsoever to the
# document's actual text, formatting, styles, or structure.
#
# There is therefore no Word object-model mutation to perform here: the
# paragraphs, runs, fonts/colors, section properties, styles and header
# content are all unchanged between "before" and "after". We simply touch
# the active document so the host sees the script ran successfully.

$d = $word.ActiveDocument
